# Week 13 logging update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append this week's per-play yardage logs to the running
# space-separated lists.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70"

# ---------------------------------------------------------------------
# OFF sheet: updated offensive totals
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 354
$offWs.Range("F2").Value = 102
$offWs.Range("G2").Value = 96
$offWs.Range("N2").Value = 38
$offWs.Range("O2").Value = 49

$offWs.Range("B3").Value = 18
$offWs.Range("C3").Value = 340
$offWs.Range("D3").Value = 11
$offWs.Range("E3").Value = 54
$offWs.Range("F3").Value = 209
$offWs.Range("G3").Value = 79
$offWs.Range("I3").Value = 106
$offWs.Range("J3").Value = 92
$offWs.Range("L3").Value = 567
$offWs.Range("M3").Value = 374
$offWs.Range("Q3").Value = 965

# ---------------------------------------------------------------------
# DEF sheet: updated defensive totals
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 367
$defWs.Range("E2").Value = 23
$defWs.Range("F2").Value = 122
$defWs.Range("G2").Value = 95
$defWs.Range("I2").Value = 13
$defWs.Range("J2").Value = 66
$defWs.Range("N2").Value = 27
$defWs.Range("O2").Value = 34
$defWs.Range("P2").Value = 18

$defWs.Range("B3").Value = 25
$defWs.Range("C3").Value = 302
$defWs.Range("F3").Value = 166
$defWs.Range("G3").Value = 68
$defWs.Range("H3").Value = 40
$defWs.Range("I3").Value = 100
$defWs.Range("J3").Value = 90
$defWs.Range("L3").Value = 460
$defWs.Range("M3").Value = 302
$defWs.Range("Q3").Value = 869

# ---------------------------------------------------------------------
# ST sheet: updated special-teams totals + logs
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 161
$stWs.Range("D2").Value = 92
$stWs.Range("F2").Value = 339
$stWs.Range("G2").Value = 327
$stWs.Range("J2").Value = 164
$stWs.Range("K2").Value = 155
$stWs.Range("N2").Value = 68
$stWs.Range("O2").Value = 38

$stWs.Range("B3").Value = 121

$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 60 59 59 54 31 37 37"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 8 19 0 8 0 2 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 7 16 9 0 0 12"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 22"

# ---------------------------------------------------------------------
# TURNS sheet: updated turnover totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 13
$turnsWs.Range("C3").Value = 14
$turnsWs.Range("E3").Value = 10

# ---------------------------------------------------------------------
# PEN sheet: updated penalty totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B3").Value = 33
